$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 997.25
$ws.Range("J17").Value = 1023.2174
$ws.Range("L17").Value = 3069.6522
$ws.Range("N17").Value = -3405.6522
$ws.Range("H19").Value = 936.4545000000001
$ws.Range("I19").Value = 799.8333
$ws.Range("J19").Value = 1100.4
$ws.Range("K19").Value = 799.8333
$ws.Range("L19").Value = 1100.4
$ws.Range("M19").Value = -624.8333
$ws.Range("N19").Value = -1450.4
$ws.Range("H33").Value = 238.8
$ws.Range("I33").Value = 111.111115
$ws.Range("K33").Value = 111.111115
$ws.Range("M33").Value = 117.888885
$ws.Range("H38").Value = 1088.3334
$ws.Range("I38").Value = 195
$ws.Range("J38").Value = 2875
$ws.Range("K38").Value = 585
$ws.Range("L38").Value = 8625
$ws.Range("M38").Value = -213
$ws.Range("N38").Value = -9369
$ws.Range("H39").Value = 127
$ws.Range("J39").Value = 240
$ws.Range("L39").Value = 720
$ws.Range("N39").Value = -1312
$ws.Range("H43").Value = 812
$ws.Range("I43").Value = 750
$ws.Range("J43").Value = 824.4
$ws.Range("K43").Value = 750
$ws.Range("L43").Value = 824.4
$ws.Range("M43").Value = -681
$ws.Range("N43").Value = -962.4
$ws.Range("H137").Value = 58825390
$ws.Range("I137").Value = 100001680
$ws.Range("J137").Value = 2129
$ws.Range("K137").Value = 300005040
$ws.Range("L137").Value = 6387
$ws.Range("M137").Value = -300002490
$ws.Range("N137").Value = -11487
$ws.Range("H138").Value = 5915976.5
$ws.Range("J138").Value = 7355330.5
$ws.Range("L138").Value = 22065991.5
$ws.Range("N138").Value = -22076271.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4937.6895
$ws.Range("I2").Value = 12567
$ws.Range("J2").Value = 922.2632
$ws.Range("K2").Value = 12567
$ws.Range("L2").Value = 922.2632
$ws.Range("M2").Value = -12454
$ws.Range("N2").Value = -1148.2632
$ws.Range("H32").Value = 2141.5144
$ws.Range("I32").Value = 1425.0938
$ws.Range("J32").Value = 9783.333000000001
$ws.Range("K32").Value = 1425.0938
$ws.Range("L32").Value = 9783.333000000001
$ws.Range("M32").Value = -1138.0938
$ws.Range("N32").Value = -10357.333
$ws.Range("H61").Value = 2055.3809
$ws.Range("I61").Value = 1484.4117
$ws.Range("J61").Value = 4482
$ws.Range("K61").Value = 1484.4117
$ws.Range("L61").Value = 4482
$ws.Range("M61").Value = -1272.4117
$ws.Range("N61").Value = -4906
$ws.Range("H112").Value = 13239.25
$ws.Range("J112").Value = 13239.25
$ws.Range("L112").Value = 13239.25
$ws.Range("N112").Value = -16193.25
$ws.Range("H116").Value = 4937.6895
$ws.Range("I116").Value = 12567
$ws.Range("J116").Value = 922.2632
$ws.Range("K116").Value = 12567
$ws.Range("L116").Value = 922.2632
$ws.Range("M116").Value = -10273
$ws.Range("N116").Value = -5510.2632
$ws.Range("H136").Value = 2055.3809
$ws.Range("I136").Value = 1484.4117
$ws.Range("J136").Value = 4482
$ws.Range("K136").Value = 4453.2351
$ws.Range("L136").Value = 13446
$ws.Range("M136").Value = -1903.2351
$ws.Range("N136").Value = -18546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4937.6895
$ws.Range("I3").Value = 12567
$ws.Range("J3").Value = 922.2632
$ws.Range("K3").Value = 12567
$ws.Range("L3").Value = 922.2632
$ws.Range("M3").Value = -12453
$ws.Range("N3").Value = -1150.2632
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1645.1305
$ws.Range("I31").Value = 1295.875
$ws.Range("J31").Value = 2443.4285
$ws.Range("K31").Value = 1295.875
$ws.Range("L31").Value = 2443.4285
$ws.Range("M31").Value = -1000.875
$ws.Range("N31").Value = -3033.4285
$ws.Range("H34").Value = 1645.1305
$ws.Range("I34").Value = 1295.875
$ws.Range("J34").Value = 2443.4285
$ws.Range("K34").Value = 1295.875
$ws.Range("L34").Value = 2443.4285
$ws.Range("M34").Value = -1093.875
$ws.Range("N34").Value = -2847.4285
$ws.Range("H58").Value = 1460.2222
$ws.Range("I58").Value = 857
$ws.Range("J58").Value = 2892.875
$ws.Range("K58").Value = 857
$ws.Range("L58").Value = 2892.875
$ws.Range("M58").Value = -654
$ws.Range("N58").Value = -3298.875
$ws.Range("H136").Value = 1460.2222
$ws.Range("I136").Value = 857
$ws.Range("J136").Value = 2892.875
$ws.Range("K136").Value = 2571
$ws.Range("L136").Value = 8678.625
$ws.Range("M136").Value = -21
$ws.Range("N136").Value = -13778.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 13158655
$ws.Range("I113").Value = 333
$ws.Range("J113").Value = 14706693
$ws.Range("K113").Value = 999
$ws.Range("L113").Value = 44120079
$ws.Range("M113").Value = 1171
$ws.Range("N113").Value = -44124419
$ws.Range("H131").Value = 3381.6667
$ws.Range("I131").Value = 276
$ws.Range("J131").Value = 3698.5715
$ws.Range("K131").Value = 828
$ws.Range("L131").Value = 11095.7145
$ws.Range("M131").Value = 4212
$ws.Range("N131").Value = -21175.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 305.55554
$ws.Range("J107").Value = 301.33334
$ws.Range("L107").Value = 301.33334
$ws.Range("N107").Value = -4141.33334
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H132").Value = 2257.5
$ws.Range("I132").Value = 1619.8823
$ws.Range("J132").Value = 3806
$ws.Range("K132").Value = 4859.6469
$ws.Range("L132").Value = 11418
$ws.Range("M132").Value = -2329.6469
$ws.Range("N132").Value = -16478

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3373.8696
$ws.Range("I7").Value = 2966.5
$ws.Range("J7").Value = 3517.647
$ws.Range("K7").Value = 2966.5
$ws.Range("L7").Value = 3517.647
$ws.Range("M7").Value = -2854.5
$ws.Range("N7").Value = -3741.647
$ws.Range("H22").Value = 21658
$ws.Range("H27").Value = 21658
$ws.Range("H46").Value = 1057.7142
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 1200.8
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 1200.8
$ws.Range("M46").Value = -512
$ws.Range("N46").Value = -1576.8
$ws.Range("H110").Value = 25857.6
$ws.Range("J110").Value = 25857.6
$ws.Range("L110").Value = 25857.6
$ws.Range("N110").Value = -34037.6
$ws.Range("H122").Value = 3600.2354
$ws.Range("I122").Value = 3040.8
$ws.Range("J122").Value = 3833.3333
$ws.Range("K122").Value = 9122.400000000001
$ws.Range("L122").Value = 11499.9999
$ws.Range("M122").Value = -6672.400000000001
$ws.Range("N122").Value = -16399.9999
$ws.Range("H126").Value = 3373.8696
$ws.Range("I126").Value = 2966.5
$ws.Range("J126").Value = 3517.647
$ws.Range("K126").Value = 8899.5
$ws.Range("L126").Value = 10552.941
$ws.Range("M126").Value = -6429.5
$ws.Range("N126").Value = -15492.941
